$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 3 to the match log, duplicating the existing row 2
# (same fixture: Dubai (DSC), RCB won by 10 runs) with Bhuvneshwar Kumar's
# figures for that innings.
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " September 21 2020"
$ws.Range("C3").Value = "RCB won by 10 runs"
$ws.Range("D3").Value = "Sunrisers Hyderabad"
$ws.Range("E3").Value = "Royal Challengers Bangalore"
$ws.Range("F3").Value = "Bhuvneshwar Kumar "

# Numeric-looking stats are stored as text (matching the rest of the sheet),
# so enter them with a leading apostrophe and then strip the resulting
# "Text" number format back to General, same as every other cell.
$ws.Range("G3").Value = "'0"
$ws.Range("H3").Value = "'2"
$ws.Range("I3").Value = "'0"
$ws.Range("J3").Value = "'0"
$ws.Range("K3").Value = "'0.00"
$ws.Range("G3:K3").ClearFormats()
